$d = $word.ActiveDocument
Write-Output "placeholder"
